# Fill in the last (previously empty) row of the Hilfsmittelverzeichnis
# table with the new ChatGPT entry, as described by the commit
# "feat: Initial analysis setup".
#
# Row 5 of the (only) table on the page is empty (three empty <w:p/>
# cells). We populate:
#   Col 1 (Tool)  : "ChatGPT, chatgpt.com"
#   Col 2 (Prompt): prompt paragraph + a second paragraph that is a
#                   hyperlink to the shared chat
#   Col 3 (Datum) : "Initialisiert" paragraph + "05.05.2025" paragraph
#
# NOTE: Cell/Range objects obtained before a text-insertion becomes
# stale once another part of the document is edited, so each cell is
# re-fetched fresh (via $d.Tables.Item(1)) immediately before it is
# used, and columns are processed from last to first (3, 2, 1) so
# that earlier character offsets in the document are never shifted
# out from under us by a later edit.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$row = $t.Rows.Count   # last row, currently empty

# --- Column 3 (Datum): two paragraphs, "Initialisiert" / "05.05.2025" ---
$t = $d.Tables.Item(1)
$cell = $t.Cell($row, 3)
$cell.Range.Text = "Initialisiert" + [char]13 + "05.05.2025"

# --- Column 2 (Prompt): prompt paragraph + hyperlink paragraph ---
$t = $d.Tables.Item(1)
$cell = $t.Cell($row, 2)
$cellStart = $cell.Range.Start

$promptText = "Context for this chat: I have a json dataset which I want to analyse with R in VSCode using the tidyverse package, which contains ggplot2 for data visualisation. The dataset file (~/analysis/dataset.json) is stored in the same folder as the r file which is used to execute the analysis. Never change values, only value formats when directly told so. Highlight the diff for each adjustment in any of the files. Document functions"
$shareUrl = "https://chatgpt.com/share/681de9bc-27f8-8007-b147-6f995028c468"

# Write both lines as plain text first (separated by a paragraph mark),
# then turn the second line into a real hyperlink run.
$cell.Range.Text = $promptText + [char]13 + $shareUrl

$urlStart = $cellStart + $promptText.Length + 1
$urlEnd = $urlStart + $shareUrl.Length
$linkRange = $d.Range($urlStart, $urlEnd)
$d.Hyperlinks.Add($linkRange, $shareUrl, [Type]::Missing, [Type]::Missing, $shareUrl)

# --- Column 1 (Tool): "ChatGPT, chatgpt.com" ---
$t = $d.Tables.Item(1)
$cell = $t.Cell($row, 1)
$cell.Range.Text = "ChatGPT, chatgpt.com"
